$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 5;   I = "sv"; J = "Statement-opinion" },
    @{ Row = 13;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 50;  I = "aa"; J = "Agree/Accept" },
    @{ Row = 54;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 55;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 98;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 99;  I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 112; I = "qy"; J = "Yes-No-Question" },
    @{ Row = 119; I = "%";  J = "Uninterpretable" },
    @{ Row = 121; I = "sv"; J = "Statement-opinion" },
    @{ Row = 133; I = "aa"; J = "Agree/Accept" },
    @{ Row = 150; I = "sv"; J = "Statement-opinion" },
    @{ Row = 153; I = "sd"; J = "Statement-non-opinion" },
    @{ Row = 166; I = "sv"; J = "Statement-opinion" },
    @{ Row = 170; I = "sv"; J = "Statement-opinion" },
    @{ Row = 176; I = "ba"; J = "Appreciation" },
    @{ Row = 184; I = "ba"; J = "Appreciation" },
    @{ Row = 191; I = "ba"; J = "Appreciation" },
    @{ Row = 194; I = "sv"; J = "Statement-opinion" },
    @{ Row = 213; I = "ba"; J = "Appreciation" },
    @{ Row = 215; I = "sd"; J = "Statement-non-opinion" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
